$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume) columns hold text-formatted values (not real numbers);
# force Text number format first so Excel does not coerce them to floats/dates.
$textCells = @('D2','E2','D3','E3','D4','E4','D5','E5','D6','E6','D7','E7','D8','E8','D9','E9','D10','E10','D11','E11','D12','E12','D13','E13','D14','E14','D15','E15','D16','E16','D17','E17','D18','E18','D19','E19','E20','D21','E21','D22','E22','E23','D24','E24','D25','E25','D26','E26','E27','D28','E28','E29','D30','E30','E31','D32','E32','D33','E33','D34','E34','D35','E35','D36','E36','D37','E37','D38','E38','D39','E39','E40','D41','E41','D42','E42','D43','E43','D44','E44','D46','E46','E47','D48','E48','E49','D50','E50','D51','E51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.285.25'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '1.676.51'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '217.58'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').Value = '0.5269'
$ws.Range('E6').Value = '  +3.81%  '
$ws.Range('D7').Value = '1.008'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '0.2687'
$ws.Range('E8').Value = '  +2.14%  '
$ws.Range('D9').Value = '0.06469'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('D10').Value = '21.88'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('D11').Value = '0.07510'
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('D12').Value = '1.693.83'
$ws.Range('E12').Value = '  +1.67%  '
$ws.Range('D13').Value = '4.510'
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('D14').Value = '0.5770'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').Value = '0.000008473'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').Value = '64.64'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').Value = '26.309.38'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = '4.915'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').Value = '1.008'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').Value = '189.90'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').Value = '6.191'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = '144.69'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = '0.1279'
$ws.Range('E25').Value = '  +7.87%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '7.801'
$ws.Range('E26').Value = '  +3.17%  '
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('D28').Value = '0.06474'
$ws.Range('E28').Value = '  -1.11%  '
$ws.Range('E29').Value = '  +4.48%  '
$ws.Range('D30').Value = '1.318'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').Value = '  +1.83%  '
$ws.Range('D32').Value = '3.583'
$ws.Range('E32').Value = '  +2.65%  '
$ws.Range('D33').Value = '1.653'
$ws.Range('E33').Value = '  +1.94%  '
$ws.Range('D34').Value = '1.028'
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').Value = '0.6180'
$ws.Range('E35').Value = '  +2.25%  '
$ws.Range('D36').Value = '2.406'
$ws.Range('E36').Value = '  +1.55%  '
$ws.Range('D37').Value = '2.736'
$ws.Range('E37').Value = '  +2.02%  '
$ws.Range('D38').Value = '6.280'
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').Value = '1.117.12'
$ws.Range('E39').Value = '  +4.01%  '
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('D41').Value = '0.8709'
$ws.Range('E41').Value = '  +1.52%  '
$ws.Range('D42').Value = '1.015'
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('D43').Value = '100.52'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').Value = '1.826.61'
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('D46').Value = '56.88'
$ws.Range('E46').Value = '  +1.33%  '
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').Value = '8.162'
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('D50').Value = '0.4293'
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('D51').Value = '6.051'
$ws.Range('E51').Value = '  +1.98%  '
